# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$wsMatches = $wb.Worksheets.Item("Matches_SOG")
$wsMatches.Cells.Item(487, 1).Value = "'897784"
$wsMatches.Cells.Item(487, 1).Style = "Normal"
$wsMatches.Cells.Item(487, 2).Value = "2025-11-18T17:00:00"
$wsMatches.Cells.Item(487, 3).Value = "Трактор"
$wsMatches.Cells.Item(487, 4).Value = "Металлург Мг"
$wsMatches.Cells.Item(487, 5).Value = 49
$wsMatches.Cells.Item(487, 6).Value = 33
$wsMatches.Cells.Item(487, 7).Value = "khl_text"

$wsMatches.Cells.Item(488, 1).Value = "'897790"
$wsMatches.Cells.Item(488, 1).Style = "Normal"
$wsMatches.Cells.Item(488, 2).Value = "2025-11-18T17:00:00"
$wsMatches.Cells.Item(488, 3).Value = "Автомобилист"
$wsMatches.Cells.Item(488, 4).Value = "Нефтехимик"
$wsMatches.Cells.Item(488, 5).Value = 33
$wsMatches.Cells.Item(488, 6).Value = 33
$wsMatches.Cells.Item(488, 7).Value = "khl_text"

$wsMatches.Cells.Item(489, 1).Value = "'897789"
$wsMatches.Cells.Item(489, 1).Style = "Normal"
$wsMatches.Cells.Item(489, 2).Value = "2025-11-18T18:00:00"
$wsMatches.Cells.Item(489, 3).Value = "Лада"
$wsMatches.Cells.Item(489, 4).Value = "Авангард"
$wsMatches.Cells.Item(489, 5).Value = 23
$wsMatches.Cells.Item(489, 6).Value = 38
$wsMatches.Cells.Item(489, 7).Value = "khl_text"

$wsMatches.Cells.Item(490, 1).Value = "'897785"
$wsMatches.Cells.Item(490, 1).Style = "Normal"
$wsMatches.Cells.Item(490, 2).Value = "2025-11-18T19:30:00"
$wsMatches.Cells.Item(490, 3).Value = "Северсталь"
$wsMatches.Cells.Item(490, 4).Value = "Динамо Мн"
$wsMatches.Cells.Item(490, 5).Value = 22
$wsMatches.Cells.Item(490, 6).Value = 20
$wsMatches.Cells.Item(490, 7).Value = "khl_text"

$wsMatches.Cells.Item(491, 1).Value = "'897786"
$wsMatches.Cells.Item(491, 1).Style = "Normal"
$wsMatches.Cells.Item(491, 2).Value = "2025-11-18T19:00:00"
$wsMatches.Cells.Item(491, 3).Value = "Ак Барс"
$wsMatches.Cells.Item(491, 4).Value = "ЦСКА"
$wsMatches.Cells.Item(491, 5).Value = 25
$wsMatches.Cells.Item(491, 6).Value = 15
$wsMatches.Cells.Item(491, 7).Value = "khl_text"

$wsMatches.Cells.Item(492, 1).Value = "'897788"
$wsMatches.Cells.Item(492, 1).Style = "Normal"
$wsMatches.Cells.Item(492, 2).Value = "2025-11-18T19:30:00"
$wsMatches.Cells.Item(492, 3).Value = "Торпедо"
$wsMatches.Cells.Item(492, 4).Value = "Адмирал"
$wsMatches.Cells.Item(492, 5).Value = 26
$wsMatches.Cells.Item(492, 6).Value = 34
$wsMatches.Cells.Item(492, 7).Value = "khl_text"

$wsMatches.Cells.Item(493, 1).Value = "'897787"
$wsMatches.Cells.Item(493, 1).Style = "Normal"
$wsMatches.Cells.Item(493, 2).Value = "2025-11-18T19:30:00"
$wsMatches.Cells.Item(493, 3).Value = "Динамо М"
$wsMatches.Cells.Item(493, 4).Value = "Спартак"
$wsMatches.Cells.Item(493, 5).Value = 28
$wsMatches.Cells.Item(493, 6).Value = 31
$wsMatches.Cells.Item(493, 7).Value = "khl_text"

$wsShotsHA = $wb.Worksheets.Item("Shots_HA")
$wsShotsHA.Range("D2").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("F2").Value = 20
$wsShotsHA.Range("K2").Value = 723
$wsShotsHA.Range("L2").Value = 616
$wsShotsHA.Range("N2").Value = 30.8
$wsShotsHA.Range("D3").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("E3").Value = 19
$wsShotsHA.Range("G3").Value = 554
$wsShotsHA.Range("H3").Value = 591
$wsShotsHA.Range("I3").Value = 29.2
$wsShotsHA.Range("J3").Value = 31.1
$wsShotsHA.Range("D4").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("F4").Value = 22
$wsShotsHA.Range("K4").Value = 698
$wsShotsHA.Range("L4").Value = 610
$wsShotsHA.Range("M4").Value = 31.7
$wsShotsHA.Range("N4").Value = 27.7
$wsShotsHA.Range("D5").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("E5").Value = 26
$wsShotsHA.Range("G5").Value = 871
$wsShotsHA.Range("H5").Value = 656
$wsShotsHA.Range("I5").Value = 33.5
$wsShotsHA.Range("J5").Value = 25.2
$wsShotsHA.Range("D6").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("D7").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("D8").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("E8").Value = 19
$wsShotsHA.Range("G8").Value = 626
$wsShotsHA.Range("H8").Value = 522
$wsShotsHA.Range("I8").Value = 32.9
$wsShotsHA.Range("J8").Value = 27.5
$wsShotsHA.Range("D9").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("F9").Value = 20
$wsShotsHA.Range("K9").Value = 715
$wsShotsHA.Range("L9").Value = 540
$wsShotsHA.Range("M9").Value = 35.8
$wsShotsHA.Range("N9").Value = 27
$wsShotsHA.Range("D10").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("D11").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("E11").Value = 24
$wsShotsHA.Range("G11").Value = 652
$wsShotsHA.Range("H11").Value = 855
$wsShotsHA.Range("I11").Value = 27.2
$wsShotsHA.Range("J11").Value = 35.6
$wsShotsHA.Range("D12").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("D13").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("F13").Value = 19
$wsShotsHA.Range("K13").Value = 556
$wsShotsHA.Range("L13").Value = 531
$wsShotsHA.Range("M13").Value = 29.3
$wsShotsHA.Range("N13").Value = 27.9
$wsShotsHA.Range("D14").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("F14").Value = 21
$wsShotsHA.Range("K14").Value = 591
$wsShotsHA.Range("L14").Value = 786
$wsShotsHA.Range("M14").Value = 28.1
$wsShotsHA.Range("N14").Value = 37.4
$wsShotsHA.Range("D15").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("D16").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("D17").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("E17").Value = 18
$wsShotsHA.Range("G17").Value = 513
$wsShotsHA.Range("H17").Value = 410
$wsShotsHA.Range("I17").Value = 28.5
$wsShotsHA.Range("J17").Value = 22.8
$wsShotsHA.Range("D18").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("D19").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("F19").Value = 16
$wsShotsHA.Range("K19").Value = 549
$wsShotsHA.Range("L19").Value = 562
$wsShotsHA.Range("M19").Value = 34.3
$wsShotsHA.Range("N19").Value = 35.1
$wsShotsHA.Range("D20").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("E20").Value = 24
$wsShotsHA.Range("G20").Value = 786
$wsShotsHA.Range("H20").Value = 727
$wsShotsHA.Range("I20").Value = 32.8
$wsShotsHA.Range("J20").Value = 30.3
$wsShotsHA.Range("D21").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("E21").Value = 20
$wsShotsHA.Range("G21").Value = 688
$wsShotsHA.Range("H21").Value = 603
$wsShotsHA.Range("I21").Value = 34.4
$wsShotsHA.Range("J21").Value = 30.1
$wsShotsHA.Range("D22").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("D23").Value = "2025-11-18T19:30:00Z"
$wsShotsHA.Range("F23").Value = 24
$wsShotsHA.Range("K23").Value = 593
$wsShotsHA.Range("L23").Value = 671
$wsShotsHA.Range("M23").Value = 24.7
$wsShotsHA.Range("N23").Value = 28

$wsShotsSummary = $wb.Worksheets.Item("Shots_Summary")
$wsShotsSummary.Range("D2").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E2").Value = 43
$wsShotsSummary.Range("F2").Value = 1474
$wsShotsSummary.Range("G2").Value = 1267
$wsShotsSummary.Range("H2").Value = 34.3
$wsShotsSummary.Range("I2").Value = 29.5
$wsShotsSummary.Range("D3").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E3").Value = 47
$wsShotsSummary.Range("F3").Value = 1335
$wsShotsSummary.Range("G3").Value = 1462
$wsShotsSummary.Range("H3").Value = 28.4
$wsShotsSummary.Range("D4").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E4").Value = 40
$wsShotsSummary.Range("F4").Value = 1384
$wsShotsSummary.Range("G4").Value = 1093
$wsShotsSummary.Range("I4").Value = 27.3
$wsShotsSummary.Range("D5").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E5").Value = 47
$wsShotsSummary.Range("F5").Value = 1582
$wsShotsSummary.Range("G5").Value = 1283
$wsShotsSummary.Range("H5").Value = 33.7
$wsShotsSummary.Range("I5").Value = 27.3
$wsShotsSummary.Range("D6").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("D7").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("D8").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E8").Value = 42
$wsShotsSummary.Range("F8").Value = 1273
$wsShotsSummary.Range("G8").Value = 1258
$wsShotsSummary.Range("H8").Value = 30.3
$wsShotsSummary.Range("I8").Value = 30
$wsShotsSummary.Range("D9").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E9").Value = 45
$wsShotsSummary.Range("F9").Value = 1635
$wsShotsSummary.Range("G9").Value = 1217
$wsShotsSummary.Range("H9").Value = 36.3
$wsShotsSummary.Range("I9").Value = 27
$wsShotsSummary.Range("D10").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("D11").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E11").Value = 45
$wsShotsSummary.Range("F11").Value = 1183
$wsShotsSummary.Range("G11").Value = 1652
$wsShotsSummary.Range("H11").Value = 26.3
$wsShotsSummary.Range("D12").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("D13").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E13").Value = 45
$wsShotsSummary.Range("F13").Value = 1482
$wsShotsSummary.Range("G13").Value = 1191
$wsShotsSummary.Range("I13").Value = 26.5
$wsShotsSummary.Range("D14").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E14").Value = 47
$wsShotsSummary.Range("F14").Value = 1399
$wsShotsSummary.Range("G14").Value = 1672
$wsShotsSummary.Range("H14").Value = 29.8
$wsShotsSummary.Range("D15").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("D16").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("D17").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E17").Value = 44
$wsShotsSummary.Range("F17").Value = 1367
$wsShotsSummary.Range("G17").Value = 1102
$wsShotsSummary.Range("H17").Value = 31.1
$wsShotsSummary.Range("I17").Value = 25
$wsShotsSummary.Range("D18").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("D19").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E19").Value = 43
$wsShotsSummary.Range("F19").Value = 1503
$wsShotsSummary.Range("G19").Value = 1308
$wsShotsSummary.Range("I19").Value = 30.4
$wsShotsSummary.Range("D20").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E20").Value = 51
$wsShotsSummary.Range("F20").Value = 1717
$wsShotsSummary.Range("G20").Value = 1593
$wsShotsSummary.Range("H20").Value = 33.7
$wsShotsSummary.Range("D21").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E21").Value = 47
$wsShotsSummary.Range("F21").Value = 1611
$wsShotsSummary.Range("G21").Value = 1471
$wsShotsSummary.Range("H21").Value = 34.3
$wsShotsSummary.Range("D22").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("D23").Value = "2025-11-18T19:30:00Z"
$wsShotsSummary.Range("E23").Value = 44
$wsShotsSummary.Range("F23").Value = 1063
$wsShotsSummary.Range("G23").Value = 1248
$wsShotsSummary.Range("H23").Value = 24.2

$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Range("B2").Value = "2025-11-18T19:30:00Z"
$wsMeta.Range("D2").Value = 74

Write-Output "done"
